$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30

$ws.Cells.Item($row, 1).Value = "Guatemala"
$ws.Cells.Item($row, 2).Value = "guatemala"
$ws.Cells.Item($row, 3).Value = "Primary Station"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "2025-10-29"
$ws.Cells.Item($row, 6).Value = 3
$ws.Cells.Item($row, 7).Value = 14.22499999999995
$ws.Cells.Item($row, 8).Value = -90.32499999999959
$ws.Cells.Item($row, 9).Value = 3
$ws.Cells.Item($row, 10).Value = 202.9380777444821
$ws.Cells.Item($row, 11).Value = "LOW"
$ws.Cells.Item($row, 12).Value = 171.182462054225
$ws.Cells.Item($row, 13).Value = 242.9454220724218
$ws.Cells.Item($row, 14).Value = 50
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 80.484375
$ws.Cells.Item($row, 18).Value = 80.58984375
$ws.Cells.Item($row, 19).Value = 79.4921875
$ws.Cells.Item($row, 20).Value = 83.46875
$ws.Cells.Item($row, 21).Value = 80.14453125
$ws.Cells.Item($row, 22).Value = 80.802734375
$ws.Cells.Item($row, 23).Value = $false
$ws.Cells.Item($row, 24).Value = -60.34042704329874
